$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 4.917747328333864
$ws.Range("D2").Value = 9.649123242017136
$ws.Range("E2").Value = 13.75295520232479
$ws.Range("F2").Value = 29.97809133218519
$ws.Range("G2").Value = 29.76393326517837
$ws.Range("H2").Value = 14.51546782401138
$ws.Range("J2").Value = 9.918164087401998
$ws.Range("K2").Value = 12.47589778731886
$ws.Range("M2").Value = 16.22776373704685
$ws.Range("N2").Value = 17.59247050624991
$ws.Range("O2").Value = 22.23613916293464
# Row 3
$ws.Range("C3").Value = 4.749897732964014
$ws.Range("D3").Value = 9.636106717478951
$ws.Range("E3").Value = 13.77058907617218
$ws.Range("F3").Value = 30.02210228545944
$ws.Range("G3").Value = 29.77942278589307
$ws.Range("H3").Value = 14.56017034592966
$ws.Range("J3").Value = 9.943897609557256
$ws.Range("K3").Value = 12.05482067005196
$ws.Range("M3").Value = 16.06776600651293
$ws.Range("N3").Value = 17.62848431112643
$ws.Range("O3").Value = 22.2978880073366
# Row 4
$ws.Range("C4").Value = 4.645211234778372
$ws.Range("D4").Value = 9.629497182937031
$ws.Range("E4").Value = 13.78368168684012
$ws.Range("F4").Value = 30.056759522099
$ws.Range("G4").Value = 29.79923387707428
$ws.Range("H4").Value = 14.5901573633214
$ws.Range("J4").Value = 9.96082688621671
$ws.Range("K4").Value = 11.78959653306535
$ws.Range("M4").Value = 15.97090735538339
$ws.Range("N4").Value = 17.65231689396271
$ws.Range("O4").Value = 22.34102999933822
# Row 5
$ws.Range("C5").Value = 4.602218754300879
$ws.Range("D5").Value = 9.627153730014427
$ws.Range("E5").Value = 13.78958679285115
$ws.Range("F5").Value = 30.0727986150435
$ws.Range("G5").Value = 29.80989077146436
$ws.Range("H5").Value = 14.60301535223498
$ws.Range("J5").Value = 9.968009884887838
$ws.Range("K5").Value = 11.67999293003123
$ws.Range("M5").Value = 15.93182083169717
$ws.Range("N5").Value = 17.66246208115841
$ws.Range("O5").Value = 22.35992192035437
# Row 6
$ws.Range("C6").Value = 4.595062209807621
$ws.Range("D6").Value = 9.626785803189316
$ws.Range("E6").Value = 13.79060174471732
$ws.Range("F6").Value = 30.07557749151018
$ws.Range("G6").Value = 29.81181616353165
$ws.Range("H6").Value = 14.6051889311958
$ws.Range("J6").Value = 9.969219790155504
$ws.Range("K6").Value = 11.66170625272201
$ws.Range("M6").Value = 15.92535481475587
$ws.Range("N6").Value = 17.66417286822892
$ws.Range("O6").Value = 22.36313799953412
# Row 7
$ws.Range("C7").Value = 4.64463266069736
$ws.Range("D7").Value = 9.629464158328782
$ws.Range("E7").Value = 13.78375901814299
$ws.Range("F7").Value = 30.05696807865281
$ws.Range("G7").Value = 29.79936714832373
$ws.Range("H7").Value = 14.59032818800737
$ws.Range("J7").Value = 9.960922607560287
$ws.Range("K7").Value = 11.7881243270998
$ws.Range("M7").Value = 15.97037861701874
$ws.Range("N7").Value = 17.65245196047401
$ws.Range("O7").Value = 22.34127947802893
# Row 8
$ws.Range("C8").Value = 4.860254629599658
$ws.Range("D8").Value = 9.644349610782227
$ws.Range("E8").Value = 13.75856529529244
$ws.Range("F8").Value = 29.99167978104645
$ws.Range("G8").Value = 29.7671329358103
$ws.Range("H8").Value = 14.53035386497014
$ws.Range("J8").Value = 9.926802928532991
$ws.Range("K8").Value = 12.33218708154996
$ws.Range("M8").Value = 16.17233235980786
$ws.Range("N8").Value = 17.60453151730681
$ws.Range("O8").Value = 22.2563431215364
# Row 9
$ws.Range("C9").Value = 5.267114222794365
$ws.Range("D9").Value = 9.684402190549445
$ws.Range("E9").Value = 13.72712627943141
$ws.Range("F9").Value = 29.92436055019945
$ws.Range("G9").Value = 29.78584544671018
$ws.Range("H9").Value = 14.43291769475627
$ws.Range("J9").Value = 9.868835636379227
$ws.Range("K9").Value = 13.33988579251507
$ws.Range("M9").Value = 16.57766705123218
$ws.Range("N9").Value = 17.52417536104154
$ws.Range("O9").Value = 22.13140834789495
# Row 10
$ws.Range("C10").Value = 5.552624964199505
$ws.Range("D10").Value = 9.720294545554895
$ws.Range("E10").Value = 13.7149638076438
$ws.Range("F10").Value = 29.91205194386992
$ws.Range("G10").Value = 29.84966453604632
$ws.Range("H10").Value = 14.3736609955584
$ws.Range("J10").Value = 9.831677083589149
$ws.Range("K10").Value = 14.03682818487487
$ws.Range("M10").Value = 16.87890664840292
$ws.Range("N10").Value = 17.47339606685334
$ws.Range("O10").Value = 22.06517828034541
# Row 11
$ws.Range("C11").Value = 5.67894749264301
$ws.Range("D11").Value = 9.737989616750198
$ws.Range("E11").Value = 13.71179935702517
$ws.Range("F11").Value = 29.91453200781976
$ws.Range("G11").Value = 29.88955105000241
$ws.Range("H11").Value = 14.3493872326394
$ws.Range("J11").Value = 9.815947423545728
$ws.Range("K11").Value = 14.34317497224652
$ws.Range("M11").Value = 17.01621987018384
$ws.Range("N11").Value = 17.45207977921299
$ws.Range("O11").Value = 22.04063266115448
# Row 12
$ws.Range("C12").Value = 5.726224535106664
$ws.Range("D12").Value = 9.74488342847658
$ws.Range("E12").Value = 13.71094086089562
$ws.Range("F12").Value = 29.91663240673037
$ws.Range("G12").Value = 29.9062105917516
$ws.Range("H12").Value = 14.34058162737077
$ws.Range("J12").Value = 9.810159505525615
$ws.Range("K12").Value = 14.4575517341962
$ws.Range("M12").Value = 17.0682171841121
$ws.Range("N12").Value = 17.44426365373607
$ws.Range("O12").Value = 22.03214300028462
# Row 13
$ws.Range("C13").Value = 5.716068110080101
$ws.Range("D13").Value = 9.743390194608363
$ws.Range("E13").Value = 13.71111065354952
$ws.Range("F13").Value = 29.91612841970725
$ws.Range("G13").Value = 29.90255359793237
$ws.Range("H13").Value = 14.34246087987215
$ws.Range("J13").Value = 9.811398544288991
$ws.Range("K13").Value = 14.43299250076118
$ws.Range("M13").Value = 17.05701927612969
$ws.Range("N13").Value = 17.44593562439559
$ws.Range("O13").Value = 22.03393555216958
# Row 14
$ws.Range("C14").Value = 5.68284844732861
$ws.Range("D14").Value = 9.738552926425303
$ws.Range("E14").Value = 13.7117219224648
$ws.Range("F14").Value = 29.91468154206675
$ws.Range("G14").Value = 29.89089050822665
$ws.Range("H14").Value = 14.34865504390738
$ws.Range("J14").Value = 9.815467871213169
$ws.Range("K14").Value = 14.35261797535916
$ws.Range("M14").Value = 17.0204979532037
$ws.Range("N14").Value = 17.45143161667765
$ws.Range("O14").Value = 22.03991806071469
# Row 15
$ws.Range("C15").Value = 5.662426403521365
$ws.Range("D15").Value = 9.735614994660487
$ws.Range("E15").Value = 13.71214057136049
$ws.Range("F15").Value = 29.91394648521394
$ws.Range("G15").Value = 29.88394888226022
$ws.Range("H15").Value = 14.35249947911403
$ws.Range("J15").Value = 9.817982395965668
$ws.Range("K15").Value = 14.30317142505617
$ws.Range("M15").Value = 16.99812637588976
$ws.Range("N15").Value = 17.45483137833432
$ws.Range("O15").Value = 22.04368745209466
# Row 16
$ws.Range("C16").Value = 5.544293759191604
$ws.Range("D16").Value = 9.719165347216711
$ws.Range("E16").Value = 13.71521820774321
$ws.Range("F16").Value = 29.91205241693719
$ws.Range("G16").Value = 29.84727586633919
$ws.Range("H16").Value = 14.37530137317372
$ws.Range("J16").Value = 9.832728655313492
$ws.Range("K16").Value = 14.01658434832757
$ws.Range("M16").Value = 16.86993505149454
$ws.Range("N16").Value = 17.47482497529791
$ws.Range("O16").Value = 22.06689494911797
# Row 17
$ws.Range("C17").Value = 5.470876814201532
$ws.Range("D17").Value = 9.709421782657737
$ws.Range("E17").Value = 13.71771233269658
$ws.Range("F17").Value = 29.9129596791701
$ws.Range("G17").Value = 29.82755515302822
$ws.Range("H17").Value = 14.38997707211652
$ws.Range("J17").Value = 9.842075505834591
$ws.Range("K17").Value = 13.8379639091931
$ws.Range("M17").Value = 16.7913362492041
$ws.Range("N17").Value = 17.48754677895933
$ws.Range("O17").Value = 22.08256370781773
# Row 18
$ws.Range("C18").Value = 5.428316902211532
$ws.Range("D18").Value = 9.703946479296768
$ws.Range("E18").Value = 13.71936985836182
$ws.Range("F18").Value = 29.91424206595592
$ws.Range("G18").Value = 29.81723451237126
$ws.Range("H18").Value = 14.39867059176237
$ws.Range("J18").Value = 9.847562079632549
$ws.Range("K18").Value = 13.73422417735951
$ws.Range("M18").Value = 16.74615682485651
$ws.Range("N18").Value = 17.49503192577197
$ws.Range("O18").Value = 22.092101313232
# Row 19
$ws.Range("C19").Value = 5.413851273987475
$ws.Range("D19").Value = 9.702114886336915
$ws.Range("E19").Value = 13.71996938407968
$ws.Range("F19").Value = 29.91480688734866
$ws.Range("G19").Value = 29.8139158140366
$ws.Range("H19").Value = 14.4016574067452
$ws.Range("J19").Value = 9.849438727581706
$ws.Range("K19").Value = 13.69893058559836
$ws.Range("M19").Value = 16.7308659797091
$ws.Range("N19").Value = 17.49759512274295
$ws.Range("O19").Value = 22.09542073271923
# Row 20
$ws.Range("C20").Value = 5.478726935860399
$ws.Range("D20").Value = 9.710445682460069
$ws.Range("E20").Value = 13.71742375716477
$ws.Range("F20").Value = 29.91278438610302
$ws.Range("G20").Value = 29.82954869631664
$ws.Range("H20").Value = 14.38838868716336
$ws.Range("J20").Value = 9.841069081825271
$ws.Range("K20").Value = 13.85708276334711
$ws.Range("M20").Value = 16.79970056593935
$ws.Range("N20").Value = 17.48617514819478
$ws.Range("O20").Value = 22.08084135134079
# Row 21
$ws.Range("C21").Value = 5.692621373875485
$ws.Range("D21").Value = 9.739968538770954
$ws.Range("E21").Value = 13.71153316220935
$ws.Range("F21").Value = 29.9150750179509
$ws.Range("G21").Value = 29.89427408327161
$ws.Range("H21").Value = 14.34682517781491
$ws.Range("J21").Value = 9.814268038477861
$ws.Range("K21").Value = 14.37627084165804
$ws.Range("M21").Value = 17.03122547464306
$ws.Range("N21").Value = 17.44981037022988
$ws.Range("O21").Value = 22.03813898125131
# Row 22
$ws.Range("C22").Value = 5.829140032493616
$ws.Range("D22").Value = 9.760387020442311
$ws.Range("E22").Value = 13.70966355975708
$ws.Range("F22").Value = 29.92333983984734
$ws.Range("G22").Value = 29.94563764077175
$ws.Range("H22").Value = 14.32191318262609
$ws.Range("J22").Value = 9.797734407629614
$ws.Range("K22").Value = 14.70605359909113
$ws.Range("M22").Value = 17.1825208581911
$ws.Range("N22").Value = 17.4275351515392
$ws.Range("O22").Value = 22.01492493753267
# Row 23
$ws.Range("C23").Value = 5.75659112541301
$ws.Range("D23").Value = 9.749387710244362
$ws.Range("E23").Value = 13.71048049301148
$ws.Range("F23").Value = 29.91830991692808
$ws.Range("G23").Value = 29.91739721558606
$ws.Range("H23").Value = 14.33500291485374
$ws.Range("J23").Value = 9.806468909401394
$ws.Range("K23").Value = 14.53094254148284
$ws.Range("M23").Value = 17.10178657951132
$ws.Range("N23").Value = 17.43928758972081
$ws.Range("O23").Value = 22.02688444371876
# Row 24
$ws.Range("C24").Value = 5.475178987197774
$ws.Range("D24").Value = 9.709982383373191
$ws.Range("E24").Value = 13.71755352555987
$ws.Range("F24").Value = 29.91286126625358
$ws.Range("G24").Value = 29.82864424692482
$ws.Range("H24").Value = 14.38910599778651
$ws.Range("J24").Value = 9.841523734285881
$ws.Range("K24").Value = 13.84844238799241
$ws.Range("M24").Value = 16.79591903011378
$ws.Range("N24").Value = 17.48679472909744
$ws.Range("O24").Value = 22.0816183797458
# Row 25
$ws.Range("C25").Value = 5.159161381112855
$ws.Range("D25").Value = 9.672419357318821
$ws.Range("E25").Value = 13.73370921650891
$ws.Range("F25").Value = 29.93605489147386
$ws.Range("G25").Value = 29.77199075773625
$ws.Range("H25").Value = 14.45711373251531
$ws.Range("J25").Value = 9.883562268589346
$ws.Range("K25").Value = 13.07442313536366
$ws.Range("M25").Value = 16.46724707371117
$ws.Range("N25").Value = 17.5444605757425
$ws.Range("O25").Value = 22.16073027034967
